# Apply BOM updates: add reference designators for the newly added
# switch/resistor rows and update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 21 (S1 - switch) previously had no Ref-Des in column A.
$ws.Range("A21").Value = "S1"

# Rows 22-24 previously all shared the placeholder Ref-Des "RX";
# give them their real, distinct resistor designators.
$ws.Range("A22").Value = "R7"
$ws.Range("A23").Value = "R9"
$ws.Range("A24").Value = "R8"

# Update the active selection to match the author's final cursor position.
$ws.Range("B24").Select()
